$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '39.555.13'
$ws.Cells.Item(2, 5).Value = '  +2.11%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.166.78'
$ws.Cells.Item(3, 5).Value = '  +3.12%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '228.94'
$ws.Cells.Item(5, 5).Value = '  +0.47%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.632'
$ws.Cells.Item(6, 5).Value = '  +2.54%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '63.79'
$ws.Cells.Item(7, 5).Value = '  +2.44%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.395'
$ws.Cells.Item(9, 5).Value = '  +1.34%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.71%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.62%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.15%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '2.488.09'
$ws.Cells.Item(13, 5).Value = '  +3.07%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '22.09'
$ws.Cells.Item(14, 5).Value = '  +0.32%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.812'
$ws.Cells.Item(15, 5).Value = '  +0.37%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '5.52'
$ws.Cells.Item(16, 5).Value = '  -0.13%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.164.39'
$ws.Cells.Item(17, 5).Value = '  +2.73%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '39.538.37'
$ws.Cells.Item(18, 5).Value = '  +2.08%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.22'
$ws.Cells.Item(19, 5).Value = '  +1.61%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '72.24'
$ws.Cells.Item(20, 5).Value = '  +0.61%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '0.0₃0848'
$ws.Cells.Item(21, 5).Value = '  +1.07%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '230.10'
$ws.Cells.Item(22, 5).Value = '  +0.97%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.05%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +1.78%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.31'
$ws.Cells.Item(25, 5).Value = '  -2.11%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '9.74'
$ws.Cells.Item(26, 5).Value = '  +1.18%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '172.14'
$ws.Cells.Item(27, 5).Value = '  -0.08%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.139'
$ws.Cells.Item(28, 5).Value = '  +0.29%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '19.90'
$ws.Cells.Item(29, 5).Value = '  +2.88%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.42'
$ws.Cells.Item(30, 5).Value = '  +0.38%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '2.66'
$ws.Cells.Item(31, 5).Value = '  +5.63%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +2.10%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '4.61'
$ws.Cells.Item(33, 5).Value = '  +1.66%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.52%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '7.05'
$ws.Cells.Item(35, 5).Value = '  +0.56%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.0621'
$ws.Cells.Item(36, 5).Value = '  +0.21%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.43'
$ws.Cells.Item(37, 5).Value = '  +1.70%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.66'
$ws.Cells.Item(38, 5).Value = '  +2.51%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.998'
$ws.Cells.Item(39, 5).Value = '  -0.17%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '17.96'
$ws.Cells.Item(40, 5).Value = '  -0.58%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '102.79'
$ws.Cells.Item(41, 5).Value = '  -0.06%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.0229'
$ws.Cells.Item(42, 5).Value = '  +0.18%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '1.526.28'
$ws.Cells.Item(43, 5).Value = '  -0.22%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +1.47%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '4.39'
$ws.Cells.Item(45, 5).Value = '  +6.20%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '7.93'
$ws.Cells.Item(46, 5).Value = '  +2.36%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'ARBITRUM'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.11'
$ws.Cells.Item(47, 5).Value = '  +5.26%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.28%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.0925'
$ws.Cells.Item(49, 5).Value = '  +1.40%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '2.371.30'
$ws.Cells.Item(50, 5).Value = '  +3.09%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.58%  '
